$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12250
$ws1.Range("F11").Value = 213
$ws1.Range("F13").Value = 5354
$ws1.Range("F15").Value = 216
$ws1.Range("F16").Value = 564

# Sheet "全部类型" (sheet4): update the same underlying rows (shifted by one
# extra row compared to "展览") with identical new values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 12250
$ws4.Range("F12").Value = 213
$ws4.Range("F14").Value = 5354
$ws4.Range("F16").Value = 216
$ws4.Range("F17").Value = 564
